# Re-order the header/column labels on row 2 of the (single) worksheet.
#
# The commit re-orders the columns of the staging template:
#   old order: BusinessKey, Code, LongName, OutputBusinessKey, ShortName, SubOutput_ID, TextDescription
#   new order: SubOutput_ID, BusinessKey, OutputBusinessKey, Code, LongName, ShortName, TextDescription
#
# Row 1 (the "For internal use only..." note) is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "SubOutput_ID"
$ws.Range("B2").Value = "BusinessKey"
$ws.Range("C2").Value = "OutputBusinessKey"
$ws.Range("D2").Value = "Code"
$ws.Range("E2").Value = "LongName"
$ws.Range("F2").Value = "ShortName"
$ws.Range("G2").Value = "TextDescription"
